# Auto-generated script to apply numeric corrections to Kujata_Profits sheets
$wb = $excel.ActiveWorkbook

# Sheet 1: ALC
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(11, 8).Value = [double]80  # H11: 79.57143000000001 -> 80
$ws.Cells.Item(11, 9).Value = [double]80  # I11: 79.57143000000001 -> 80
$ws.Cells.Item(11, 11).Value = [double]80  # K11: 79.57143000000001 -> 80
$ws.Cells.Item(11, 13).Value = [double]60  # M11: 60.42856999999999 -> 60
$ws.Cells.Item(112, 8).Value = [double]1913.1555  # H112: 2039.1025 -> 1913.1555
$ws.Cells.Item(112, 9).Value = [double]914.2857  # I112: 937.5 -> 914.2857
$ws.Cells.Item(112, 10).Value = [double]2097.158  # J112: 2323.3872 -> 2097.158
$ws.Cells.Item(112, 11).Value = [double]2742.8571  # K112: 2812.5 -> 2742.8571
$ws.Cells.Item(112, 12).Value = [double]6291.474  # L112: 6970.1616 -> 6291.474
$ws.Cells.Item(112, 13).Value = [double]-1634.8571  # M112: -1704.5 -> -1634.8571
$ws.Cells.Item(112, 14).Value = [double]-8507.474  # N112: -9186.161599999999 -> -8507.474
$ws.Cells.Item(137, 8).Value = [double]1148.0156  # H137: 1160.6984 -> 1148.0156
$ws.Cells.Item(137, 9).Value = [double]832.4865  # I137: 845.9167 -> 832.4865
$ws.Cells.Item(137, 11).Value = [double]2497.4595  # K137: 2537.7501 -> 2497.4595
$ws.Cells.Item(137, 13).Value = [double]52.54050000000007  # M137: 12.2498999999998 -> 52.54050000000007
$ws.Cells.Item(138, 8).Value = [double]1397.1771  # H138: 1400.2783 -> 1397.1771
$ws.Cells.Item(138, 10).Value = [double]1634.3549  # J138: 1635.3651 -> 1634.3549
$ws.Cells.Item(138, 12).Value = [double]4903.0647  # L138: 4906.0953 -> 4903.0647
$ws.Cells.Item(138, 14).Value = [double]-15183.0647  # N138: -15186.0953 -> -15183.0647
# Sheet 2: ARM
$ws = $wb.Worksheets.Item(2)
$ws.Cells.Item(32, 8).Value = [double]3760.9753  # H32: 4018.3333 -> 3760.9753
$ws.Cells.Item(32, 9).Value = [double]3499.697  # I32: 3744.5247 -> 3499.697
$ws.Cells.Item(32, 10).Value = [double]4910.6  # J32: 5211.357 -> 4910.6
$ws.Cells.Item(32, 11).Value = [double]3499.697  # K32: 3744.5247 -> 3499.697
$ws.Cells.Item(32, 12).Value = [double]4910.6  # L32: 5211.357 -> 4910.6
$ws.Cells.Item(32, 13).Value = [double]-3212.697  # M32: -3457.5247 -> -3212.697
$ws.Cells.Item(32, 14).Value = [double]-5484.6  # N32: -5785.357 -> -5484.6
$ws.Cells.Item(61, 8).Value = [double]17544796  # H61: 19231774 -> 17544796
$ws.Cells.Item(61, 9).Value = [double]24390864  # I61: 27027690 -> 24390864
$ws.Cells.Item(61, 10).Value = [double]1747  # J61: 1843.4667 -> 1747
$ws.Cells.Item(61, 11).Value = [double]24390864  # K61: 27027690 -> 24390864
$ws.Cells.Item(61, 12).Value = [double]1747  # L61: 1843.4667 -> 1747
$ws.Cells.Item(61, 13).Value = [double]-24390652  # M61: -27027478 -> -24390652
$ws.Cells.Item(61, 14).Value = [double]-2171  # N61: -2267.4667 -> -2171
$ws.Cells.Item(74, 8).Value = [double]1504.9166  # H74: 1547.125 -> 1504.9166
$ws.Cells.Item(74, 9).Value = [double]1089.1111  # I74: 1144.6111 -> 1089.1111
$ws.Cells.Item(74, 10).Value = [double]2752.3333  # J74: 2754.6667 -> 2752.3333
$ws.Cells.Item(74, 11).Value = [double]1089.1111  # K74: 1144.6111 -> 1089.1111
$ws.Cells.Item(74, 12).Value = [double]2752.3333  # L74: 2754.6667 -> 2752.3333
$ws.Cells.Item(74, 13).Value = [double]-215.1111000000001  # M74: -270.6111000000001 -> -215.1111000000001
$ws.Cells.Item(74, 14).Value = [double]-4500.3333  # N74: -4502.6667 -> -4500.3333
$ws.Cells.Item(77, 8).Value = [double]1504.9166  # H77: 1547.125 -> 1504.9166
$ws.Cells.Item(77, 9).Value = [double]1089.1111  # I77: 1144.6111 -> 1089.1111
$ws.Cells.Item(77, 10).Value = [double]2752.3333  # J77: 2754.6667 -> 2752.3333
$ws.Cells.Item(77, 11).Value = [double]5445.5555  # K77: 5723.0555 -> 5445.5555
$ws.Cells.Item(77, 12).Value = [double]13761.6665  # L77: 13773.3335 -> 13761.6665
$ws.Cells.Item(77, 13).Value = [double]-1077.5555  # M77: -1355.0555 -> -1077.5555
$ws.Cells.Item(77, 14).Value = [double]-22497.6665  # N77: -22509.3335 -> -22497.6665
$ws.Cells.Item(136, 8).Value = [double]17544796  # H136: 19231774 -> 17544796
$ws.Cells.Item(136, 9).Value = [double]24390864  # I136: 27027690 -> 24390864
$ws.Cells.Item(136, 10).Value = [double]1747  # J136: 1843.4667 -> 1747
$ws.Cells.Item(136, 11).Value = [double]73172592  # K136: 81083070 -> 73172592
$ws.Cells.Item(136, 12).Value = [double]5241  # L136: 5530.4001 -> 5241
$ws.Cells.Item(136, 13).Value = [double]-73170042  # M136: -81080520 -> -73170042
$ws.Cells.Item(136, 14).Value = [double]-10341  # N136: -10630.4001 -> -10341
$ws.Cells.Item(139, 8).Value = [double]54430  # H139: 75000 -> 54430
$ws.Cells.Item(139, 10).Value = [double]54430  # J139: 75000 -> 54430
$ws.Cells.Item(139, 12).Value = [double]54430  # L139: 75000 -> 54430
$ws.Cells.Item(139, 14).Value = [double]-64710  # N139: -85280 -> -64710
# Sheet 3: BSM
$ws = $wb.Worksheets.Item(3)
$ws.Cells.Item(134, 8).Value = [double]3925  # H134: 4025.342 -> 3925
$ws.Cells.Item(134, 9).Value = [double]946.5278  # I134: 970.3714 -> 946.5278
$ws.Cells.Item(134, 11).Value = [double]2839.5834  # K134: 2911.1142 -> 2839.5834
$ws.Cells.Item(134, 13).Value = [double]-304.5834  # M134: -376.1142 -> -304.5834
# Sheet 4: CRP
$ws = $wb.Worksheets.Item(4)
$ws.Cells.Item(16, 8).Value = [double]200001800  # H16: 66668096 -> 200001800
$ws.Cells.Item(16, 9).Value = [double]250001490  # I16: 66668096 -> 250001490
$ws.Cells.Item(16, 10).Value = [double]3000  # J16: 0 -> 3000
$ws.Cells.Item(16, 11).Value = [double]250001490  # K16: 66668096 -> 250001490
$ws.Cells.Item(16, 12).Value = [double]3000  # L16: 0 -> 3000
$ws.Cells.Item(16, 13).Value = [double]-250001203  # M16: -66667809 -> -250001203
$ws.Cells.Item(16, 14).Value = [double]-3574  # N16: None -> -3574
$ws.Cells.Item(31, 8).Value = [double]1802.2444  # H31: 1942.275 -> 1802.2444
$ws.Cells.Item(31, 9).Value = [double]1630.8462  # I31: 1770.3823 -> 1630.8462
$ws.Cells.Item(31, 11).Value = [double]1630.8462  # K31: 1770.3823 -> 1630.8462
$ws.Cells.Item(31, 13).Value = [double]-1335.8462  # M31: -1475.3823 -> -1335.8462
$ws.Cells.Item(34, 8).Value = [double]1802.2444  # H34: 1942.275 -> 1802.2444
$ws.Cells.Item(34, 9).Value = [double]1630.8462  # I34: 1770.3823 -> 1630.8462
$ws.Cells.Item(34, 11).Value = [double]1630.8462  # K34: 1770.3823 -> 1630.8462
$ws.Cells.Item(34, 13).Value = [double]-1428.8462  # M34: -1568.3823 -> -1428.8462
$ws.Cells.Item(56, 8).Value = [double]0  # H56: 2980 -> 0
$ws.Cells.Item(56, 9).Value = [double]0  # I56: 2980 -> 0
$ws.Cells.Item(56, 11).Value = [double]0  # K56: 2980 -> 0
$ws.Cells.Item(56, 13).ClearContents()  # M56: was -2135
$ws.Cells.Item(58, 8).Value = [double]726.63635  # H58: 765.87805 -> 726.63635
$ws.Cells.Item(58, 9).Value = [double]642.1539  # I58: 679.80554 -> 642.1539
$ws.Cells.Item(58, 11).Value = [double]642.1539  # K58: 679.80554 -> 642.1539
$ws.Cells.Item(58, 13).Value = [double]-439.1539  # M58: -476.80554 -> -439.1539
$ws.Cells.Item(107, 8).Value = [double]562.3333  # H107: 676.4375 -> 562.3333
$ws.Cells.Item(107, 9).Value = [double]396.7  # I107: 595.5714 -> 396.7
$ws.Cells.Item(107, 10).Value = [double]769.375  # J107: 739.3333 -> 769.375
$ws.Cells.Item(107, 11).Value = [double]396.7  # K107: 595.5714 -> 396.7
$ws.Cells.Item(107, 12).Value = [double]769.375  # L107: 739.3333 -> 769.375
$ws.Cells.Item(107, 13).Value = [double]1523.3  # M107: 1324.4286 -> 1523.3
$ws.Cells.Item(107, 14).Value = [double]-4609.375  # N107: -4579.3333 -> -4609.375
$ws.Cells.Item(113, 8).Value = [double]200001800  # H113: 66668096 -> 200001800
$ws.Cells.Item(113, 9).Value = [double]250001490  # I113: 66668096 -> 250001490
$ws.Cells.Item(113, 10).Value = [double]3000  # J113: 0 -> 3000
$ws.Cells.Item(113, 11).Value = [double]250001490  # K113: 66668096 -> 250001490
$ws.Cells.Item(113, 12).Value = [double]3000  # L113: 0 -> 3000
$ws.Cells.Item(113, 13).Value = [double]-249999320  # M113: -66665926 -> -249999320
$ws.Cells.Item(113, 14).Value = [double]-7340  # N113: None -> -7340
$ws.Cells.Item(132, 8).Value = [double]3420.698  # H132: 3675.898 -> 3420.698
$ws.Cells.Item(132, 9).Value = [double]3859.075  # I132: 4153.4863 -> 3859.075
$ws.Cells.Item(132, 10).Value = [double]2071.8462  # J132: 2203.3333 -> 2071.8462
$ws.Cells.Item(132, 11).Value = [double]11577.225  # K132: 12460.4589 -> 11577.225
$ws.Cells.Item(132, 12).Value = [double]6215.5386  # L132: 6609.999899999999 -> 6215.5386
$ws.Cells.Item(132, 13).Value = [double]-9047.224999999999  # M132: -9930.458899999998 -> -9047.224999999999
$ws.Cells.Item(132, 14).Value = [double]-11275.5386  # N132: -11669.9999 -> -11275.5386
$ws.Cells.Item(134, 8).Value = [double]11112441  # H134: 11629289 -> 11112441
$ws.Cells.Item(134, 9).Value = [double]1522.6  # I134: 1560.931 -> 1522.6
$ws.Cells.Item(134, 10).Value = [double]33334276  # J134: 35715296 -> 33334276
$ws.Cells.Item(134, 11).Value = [double]4567.799999999999  # K134: 4682.793 -> 4567.799999999999
$ws.Cells.Item(134, 12).Value = [double]100002828  # L134: 107145888 -> 100002828
$ws.Cells.Item(134, 13).Value = [double]-2032.799999999999  # M134: -2147.793 -> -2032.799999999999
$ws.Cells.Item(134, 14).Value = [double]-100007898  # N134: -107150958 -> -100007898
$ws.Cells.Item(136, 8).Value = [double]726.63635  # H136: 765.87805 -> 726.63635
$ws.Cells.Item(136, 9).Value = [double]642.1539  # I136: 679.80554 -> 642.1539
$ws.Cells.Item(136, 11).Value = [double]1926.4617  # K136: 2039.41662 -> 1926.4617
$ws.Cells.Item(136, 13).Value = [double]623.5382999999999  # M136: 510.58338 -> 623.5382999999999
# Sheet 5: CUL
$ws = $wb.Worksheets.Item(5)
$ws.Cells.Item(39, 8).Value = [double]2585.5  # H39: 2805.7827 -> 2585.5
$ws.Cells.Item(39, 10).Value = [double]2491.9546  # J39: 2743.842 -> 2491.9546
$ws.Cells.Item(39, 12).Value = [double]7475.8638  # L39: 8231.526 -> 7475.8638
$ws.Cells.Item(39, 14).Value = [double]-8063.8638  # N39: -8819.526 -> -8063.8638
$ws.Cells.Item(41, 8).Value = [double]1109.2  # H41: 489.73334 -> 1109.2
$ws.Cells.Item(41, 10).Value = [double]2250  # J41: 525 -> 2250
$ws.Cells.Item(41, 12).Value = [double]6750  # L41: 1575 -> 6750
$ws.Cells.Item(41, 14).Value = [double]-7426  # N41: -2251 -> -7426
$ws.Cells.Item(46, 8).Value = [double]881  # H46: 1067.5 -> 881
$ws.Cells.Item(46, 9).Value = [double]405  # I46: 1202.5 -> 405
$ws.Cells.Item(46, 11).Value = [double]1215  # K46: 3607.5 -> 1215
$ws.Cells.Item(46, 13).Value = [double]-1124  # M46: -3516.5 -> -1124
$ws.Cells.Item(55, 8).Value = [double]2413.125  # H55: 2415 -> 2413.125
$ws.Cells.Item(55, 10).Value = [double]3050.8333  # J55: 3181 -> 3050.8333
$ws.Cells.Item(55, 12).Value = [double]9152.499899999999  # L55: 9543 -> 9152.499899999999
$ws.Cells.Item(55, 14).Value = [double]-9506.499899999999  # N55: -9897 -> -9506.499899999999
$ws.Cells.Item(107, 8).Value = [double]5637.6313  # H107: 7099.6665 -> 5637.6313
$ws.Cells.Item(107, 9).Value = [double]385.33334  # I107: 442.4 -> 385.33334
$ws.Cells.Item(107, 10).Value = [double]8061.769  # J107: 10428.3 -> 8061.769
$ws.Cells.Item(107, 11).Value = [double]1156.00002  # K107: 1327.2 -> 1156.00002
$ws.Cells.Item(107, 12).Value = [double]24185.307  # L107: 31284.9 -> 24185.307
$ws.Cells.Item(107, 13).Value = [double]763.9999800000001  # M107: 592.8000000000002 -> 763.9999800000001
$ws.Cells.Item(107, 14).Value = [double]-28025.307  # N107: -35124.89999999999 -> -28025.307
$ws.Cells.Item(121, 8).Value = [double]679  # H121: 428.4 -> 679
$ws.Cells.Item(121, 9).Value = [double]425  # I121: 277.25 -> 425
$ws.Cells.Item(121, 10).Value = [double]933  # J121: 1033 -> 933
$ws.Cells.Item(121, 11).Value = [double]1275  # K121: 831.75 -> 1275
$ws.Cells.Item(121, 12).Value = [double]2799  # L121: 3099 -> 2799
$ws.Cells.Item(121, 13).Value = [double]35  # M121: 478.25 -> 35
$ws.Cells.Item(121, 14).Value = [double]-5419  # N121: -5719 -> -5419
$ws.Cells.Item(122, 8).Value = [double]859.53845  # H122: 793.4375 -> 859.53845
$ws.Cells.Item(122, 9).Value = [double]633.63635  # I122: 610.3333 -> 633.63635
$ws.Cells.Item(122, 10).Value = [double]1025.2  # J122: 903.3 -> 1025.2
$ws.Cells.Item(122, 11).Value = [double]5702.72715  # K122: 5492.9997 -> 5702.72715
$ws.Cells.Item(122, 12).Value = [double]9226.800000000001  # L122: 8129.7 -> 9226.800000000001
$ws.Cells.Item(122, 13).Value = [double]-3252.72715  # M122: -3042.9997 -> -3252.72715
$ws.Cells.Item(122, 14).Value = [double]-14126.8  # N122: -13029.7 -> -14126.8
$ws.Cells.Item(131, 8).Value = [double]21740410  # H131: 14286831 -> 21740410
$ws.Cells.Item(131, 10).Value = [double]1481  # J131: 1209.6167 -> 1481
$ws.Cells.Item(131, 12).Value = [double]4443  # L131: 3628.8501 -> 4443
$ws.Cells.Item(131, 14).Value = [double]-14523  # N131: -13708.8501 -> -14523
$ws.Cells.Item(132, 8).Value = [double]1373.7391  # H132: 1528.1904 -> 1373.7391
$ws.Cells.Item(132, 9).Value = [double]1010.2143  # I132: 1010.2857 -> 1010.2143
$ws.Cells.Item(132, 10).Value = [double]1939.2222  # J132: 2564 -> 1939.2222
$ws.Cells.Item(132, 11).Value = [double]9091.9287  # K132: 9092.5713 -> 9091.9287
$ws.Cells.Item(132, 12).Value = [double]17452.9998  # L132: 23076 -> 17452.9998
$ws.Cells.Item(132, 13).Value = [double]-6561.9287  # M132: -6562.5713 -> -6561.9287
$ws.Cells.Item(132, 14).Value = [double]-22512.9998  # N132: -28136 -> -22512.9998
# Sheet 6: GSM
$ws = $wb.Worksheets.Item(6)
$ws.Cells.Item(86, 8).Value = [double]30983.285  # H86: 28987.223 -> 30983.285
$ws.Cells.Item(86, 10).Value = [double]30983.285  # J86: 28987.223 -> 30983.285
$ws.Cells.Item(86, 12).Value = [double]30983.285  # L86: 28987.223 -> 30983.285
$ws.Cells.Item(86, 14).Value = [double]-33355.285  # N86: -31359.223 -> -33355.285
$ws.Cells.Item(89, 8).Value = [double]30983.285  # H89: 28987.223 -> 30983.285
$ws.Cells.Item(89, 10).Value = [double]30983.285  # J89: 28987.223 -> 30983.285
$ws.Cells.Item(89, 12).Value = [double]92949.855  # L89: 86961.66900000001 -> 92949.855
$ws.Cells.Item(89, 14).Value = [double]-104805.855  # N89: -98817.66900000001 -> -104805.855
$ws.Cells.Item(102, 8).Value = [double]19231668  # H102: 20834240 -> 19231668
$ws.Cells.Item(102, 9).Value = [double]31250832  # I102: 35715124 -> 31250832
$ws.Cells.Item(102, 11).Value = [double]31250832  # K102: 35715124 -> 31250832
$ws.Cells.Item(102, 13).Value = [double]-31249210  # M102: -35713502 -> -31249210
$ws.Cells.Item(113, 8).Value = [double]2000  # H113: 1612.1 -> 2000
$ws.Cells.Item(113, 9).Value = [double]0  # I113: 1560.1428 -> 0
$ws.Cells.Item(113, 10).Value = [double]2000  # J113: 1733.3334 -> 2000
$ws.Cells.Item(113, 11).Value = [double]0  # K113: 1560.1428 -> 0
$ws.Cells.Item(113, 12).Value = [double]2000  # L113: 1733.3334 -> 2000
$ws.Cells.Item(113, 13).ClearContents()  # M113: was 609.8571999999999
$ws.Cells.Item(113, 14).Value = [double]-6340  # N113: -6073.3334 -> -6340
$ws.Cells.Item(132, 8).Value = [double]1861.2858  # H132: 2726.0344 -> 1861.2858
$ws.Cells.Item(132, 9).Value = [double]1401.6786  # I132: 2416.2273 -> 1401.6786
$ws.Cells.Item(132, 11).Value = [double]4205.0358  # K132: 7248.6819 -> 4205.0358
$ws.Cells.Item(132, 13).Value = [double]-1675.0358  # M132: -4718.6819 -> -1675.0358
# Sheet 7: LTW
$ws = $wb.Worksheets.Item(7)
$ws.Cells.Item(22, 8).Value = [double]1750.5  # H22: 1499.3334 -> 1750.5
$ws.Cells.Item(22, 9).Value = [double]0  # I22: 498 -> 0
$ws.Cells.Item(22, 10).Value = [double]1750.5  # J22: 2000 -> 1750.5
$ws.Cells.Item(22, 11).Value = [double]0  # K22: 498 -> 0
$ws.Cells.Item(22, 12).Value = [double]1750.5  # L22: 2000 -> 1750.5
$ws.Cells.Item(22, 13).ClearContents()  # M22: was -203
$ws.Cells.Item(22, 14).Value = [double]-2340.5  # N22: -2590 -> -2340.5
$ws.Cells.Item(27, 8).Value = [double]1750.5  # H27: 1499.3334 -> 1750.5
$ws.Cells.Item(27, 9).Value = [double]0  # I27: 498 -> 0
$ws.Cells.Item(27, 10).Value = [double]1750.5  # J27: 2000 -> 1750.5
$ws.Cells.Item(27, 11).Value = [double]0  # K27: 498 -> 0
$ws.Cells.Item(27, 12).Value = [double]1750.5  # L27: 2000 -> 1750.5
$ws.Cells.Item(27, 13).ClearContents()  # M27: was -391
$ws.Cells.Item(27, 14).Value = [double]-1964.5  # N27: -2214 -> -1964.5
$ws.Cells.Item(100, 8).Value = [double]1725.125  # H100: 1999.8 -> 1725.125
$ws.Cells.Item(100, 9).Value = [double]1466.8334  # I100: 1666.3334 -> 1466.8334
$ws.Cells.Item(100, 11).Value = [double]1466.8334  # K100: 1666.3334 -> 1466.8334
$ws.Cells.Item(100, 13).Value = [double]-925.8334  # M100: -1125.3334 -> -925.8334
$ws.Cells.Item(123, 8).Value = [double]0  # H123: 40930 -> 0
$ws.Cells.Item(123, 10).Value = [double]0  # J123: 40930 -> 0
$ws.Cells.Item(123, 12).Value = [double]0  # L123: 40930 -> 0
$ws.Cells.Item(123, 14).ClearContents()  # N123: was -50730
$ws.Cells.Item(136, 8).Value = [double]5782.2  # H136: 4998.6 -> 5782.2
$ws.Cells.Item(136, 9).Value = [double]7264.1665  # I136: 6143.591 -> 7264.1665
$ws.Cells.Item(136, 10).Value = [double]1971.4286  # J136: 1849.875 -> 1971.4286
$ws.Cells.Item(136, 11).Value = [double]21792.4995  # K136: 18430.773 -> 21792.4995
$ws.Cells.Item(136, 12).Value = [double]5914.2858  # L136: 5549.625 -> 5914.2858
$ws.Cells.Item(136, 13).Value = [double]-19242.4995  # M136: -15880.773 -> -19242.4995
$ws.Cells.Item(136, 14).Value = [double]-11014.2858  # N136: -10649.625 -> -11014.2858
# Sheet 8: WVR
$ws = $wb.Worksheets.Item(8)
$ws.Cells.Item(109, 8).Value = [double]40359.5  # H109: 37560.5 -> 40359.5
$ws.Cells.Item(109, 10).Value = [double]30377  # J109: 33300 -> 30377
$ws.Cells.Item(109, 12).Value = [double]30377  # L109: 33300 -> 30377
$ws.Cells.Item(109, 14).Value = [double]-33151  # N109: -36074 -> -33151
$ws.Cells.Item(132, 8).Value = [double]1595.079  # H132: 1549.4103 -> 1595.079
$ws.Cells.Item(132, 9).Value = [double]1682.7059  # I132: 1524.4736 -> 1682.7059
$ws.Cells.Item(132, 10).Value = [double]1524.1428  # J132: 1573.1 -> 1524.1428
$ws.Cells.Item(132, 11).Value = [double]5048.1177  # K132: 4573.4208 -> 5048.1177
$ws.Cells.Item(132, 12).Value = [double]4572.428400000001  # L132: 4719.299999999999 -> 4572.428400000001
$ws.Cells.Item(132, 13).Value = [double]-2518.1177  # M132: -2043.4208 -> -2518.1177
$ws.Cells.Item(132, 14).Value = [double]-9632.428400000001  # N132: -9779.299999999999 -> -9632.428400000001
$ws.Cells.Item(136, 8).Value = [double]1078.2858  # H136: 780.4 -> 1078.2858
$ws.Cells.Item(136, 9).Value = [double]988.6  # I136: 672.8182 -> 988.6
$ws.Cells.Item(136, 10).Value = [double]1302.5  # J136: 1076.25 -> 1302.5
$ws.Cells.Item(136, 11).Value = [double]2965.8  # K136: 2018.4546 -> 2965.8
$ws.Cells.Item(136, 12).Value = [double]3907.5  # L136: 3228.75 -> 3907.5
$ws.Cells.Item(136, 13).Value = [double]-415.8000000000002  # M136: 531.5454 -> -415.8000000000002
$ws.Cells.Item(136, 14).Value = [double]-9007.5  # N136: -8328.75 -> -9007.5
